$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Java/Python" language notes on a few existing question rows.
# C3 gets a new (non-bold) style - vertical-top alignment, default font.
$c3 = $ws.Range("C3")
$c3.Value = "Java/Python"
$c3.VerticalAlignment = -4160
$c3.Font.Bold = $false

$ws.Range("C8").Value = "Java/Python"
$ws.Range("C15").Value = "Java/Python"

# New question row appended at the bottom of the table.
$ws.Range("A29").Value = "GFG"
$ws.Range("B29").Value = "Find length of Loop"
$ws.Range("C29").Value = "Java"

# Update the view state to match where the author was scrolled to / had
# selected when they saved.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F26").Select()
